# Deal Chat Testcases Updated
#
# 1. "Add Deal" sheet: three new Deal-Chat rows appended (28-30), cloned from
#    the existing "AutoDeal / Cottonseed / Miles (MI) / Auger" template row
#    (row 2), with new Automation Test IDs and incrementing Rate / NO of Loads.
# 2. "Deal Widget" sheet: three matching rows appended (57-59), cloned from the
#    existing "AutoDeal / Draft / SHARE" template row (row 49).
# 3. View-state bookkeeping: selection/active-cell housekeeping on the two
#    edited sheets, and the workbook's active tab moving back to "Add Deal".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Add Deal" sheet (sheet1) - append rows 28, 29, 30
# ---------------------------------------------------------------------------
$addDeal = $wb.Worksheets.Item("Add Deal")

# Row 2 is the template: AutoDeal / Cottonseed / Current Date / Current Date /
# Miles (MI) / Auger / San Francisco.../San Jose... / Added new Deal successfully
$addDeal.Range("A2:L2").Copy($addDeal.Range("A28:L28"))
$addDeal.Range("A2:L2").Copy($addDeal.Range("A29:L29"))
$addDeal.Range("A2:L2").Copy($addDeal.Range("A30:L30"))

$addDeal.Range("A28").Value = "Deals_Chat_ShipperUser_TC001"
$addDeal.Range("A29").Value = "Deals_Chat_ShipperAdmin_TC002"
$addDeal.Range("A30").Value = "Deals_Chat_CarrierUser_TC003"

$addDeal.Range("F28").Value = 10
$addDeal.Range("F29").Value = 11
$addDeal.Range("F30").Value = 12

$addDeal.Range("I28").Value = 1
$addDeal.Range("I29").Value = 2
$addDeal.Range("I30").Value = 3

# ---------------------------------------------------------------------------
# "Deal Widget" sheet (sheet3) - append rows 57, 58, 59
# ---------------------------------------------------------------------------
$dealWidget = $wb.Worksheets.Item("Deal Widget")

# Row 49 is the template: AutoDeal / Draft / SHARE / widget handled successfully
$dealWidget.Range("A49:E49").Copy($dealWidget.Range("A57:E57"))
$dealWidget.Range("A49:E49").Copy($dealWidget.Range("A58:E58"))
$dealWidget.Range("A49:E49").Copy($dealWidget.Range("A59:E59"))

$dealWidget.Range("A57").Value = "Deals_Chat_ShipperUser_TC001"
$dealWidget.Range("A58").Value = "Deals_Chat_ShipperAdmin_TC002"
$dealWidget.Range("A59").Value = "Deals_Chat_CarrierUser_TC003"

# ---------------------------------------------------------------------------
# View-state: selections on the touched sheets + "Counter Deal" no longer tab
# ---------------------------------------------------------------------------
$counterDeal = $wb.Worksheets.Item("Counter Deal")

$dealWidget.Range("C67").Select()
$addDeal.Range("C32").Select()

# Active tab moves from "Counter Deal" back to "Add Deal".
$addDeal.Activate()
$addDeal.Range("C32").Select()
